{"js": "// Apply the \"Case Study Campus 2021\" requirements-doc edits.\n// Only genuine text-level changes are applied here; the many run-splitting /\n// <w:proofErr> clean-ups in the source diff collapse to the same visible\n// text and are not observable through the Word content APIs, so they are\n// skipped.\n\nasync function replaceOnce(searchText, replacement, matchCase) {\n  const results = context.document.body.search(searchText, {\n    matchCase: matchCase === undefined ? true : matchCase,\n    matchWholeWord: false\n  });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// A. \"Birthdays\" -> \"Birthday\" (bullet under \"The Database should have\n//    following details stored\")\nawait replaceOnce(\"Birthdays\", \"Birthday\");\n\n// B. Collapse the double space before \"by\".\nawait replaceOnce(\"stored  by\", \"stored by\");\n\n// C. Insert the new \"NameNotFoundException\" bullet right after the\n//    \"...entering employees  name.\" bullet and before the\n//    \"Should be able edit ...\" bullet (same numbered list).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"entering employees\") !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchorParagraph) {\n  throw new Error(\"Could not find anchor paragraph for the new bullet.\");\n}\nanchorParagraph.insertParagraph(\n  \"If the employee with a given name is not Found Application Should throw a user defined exception NameNotFoundException.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// D. Collapse the double space between \"able\" and \"edit\".\nawait replaceOnce(\"able  edit\", \"able edit\");\n\n// E. \"Should be to get\" -> \"Should be able to get\"\nawait replaceOnce(\"Should be to get\", \"Should be able to get\");\n\n// F. Drop \"are the \" from \"The Following are the expectations\".\nawait replaceOnce(\"The Following are the expectations\", \"The Following expectations\");\n\n// G. \"4 parameter\" -> \"5 parameter\" (evaluation criteria count).\nawait replaceOnce(\"will be based on 4 parameter\", \"will be based on 5 parameter\");\n", "ps1": "# Apply the \"Case Study Campus 2021\" requirements-doc edits.\n# Only genuine text-level changes are applied here; the many run-splitting /\n# proofing-error clean-ups in the source diff collapse to the same visible\n# text and are not observable through the Word object model, so they are\n# skipped.\n\nfunction Replace-Once($doc, $findText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop - don't wrap, fail if not found\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$replaceText, 2) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\n# A. \"Birthdays\" -> \"Birthday\" (bullet under \"The Database should have\n#    following details stored\")\nReplace-Once $d \"Birthdays\" \"Birthday\"\n\n# B. Collapse the double space before \"by\".\nReplace-Once $d \"stored  by\" \"stored by\"\n\n# C. Insert the new \"NameNotFoundException\" bullet right after the\n#    \"...entering employees  name.\" bullet and before the\n#    \"Should be able edit ...\" bullet (same numbered list).\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*entering employees*\") {\n        $anchor = $d.Paragraphs.Item($i).Range\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Could not find anchor paragraph for the new bullet.\"\n}\n$anchor.Collapse(0)  # wdCollapseEnd\n$anchor.InsertAfter(\"If the employee with a given name is not Found Application Should throw a user defined exception NameNotFoundException.`r\")\n\n# D. Collapse the double space between \"able\" and \"edit\".\nReplace-Once $d \"able  edit\" \"able edit\"\n\n# E. \"Should be to get\" -> \"Should be able to get\"\nReplace-Once $d \"Should be to get\" \"Should be able to get\"\n\n# F. Drop \"are the \" from \"The Following are the expectations\".\nReplace-Once $d \"The Following are the expectations\" \"The Following expectations\"\n\n# G. \"4 parameter\" -> \"5 parameter\" (evaluation criteria count).\nReplace-Once $d \"will be based on 4 parameter\" \"will be based on 5 parameter\"\n"}
